$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.385.36'
$ws.Range("E2").Value = '  +0.49%  '

$ws.Range("D3").Value = '2.013.76'
$ws.Range("E3").Value = '  -0.41%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '260.47'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.30%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.621'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.20%  '

$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '56.87'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -5.39%  '

$ws.Range("E9").Value = '  -1.87%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0775'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.49%  '

$ws.Range("E11").Value = '  -2.75%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.32'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -5.88%  '

$ws.Range("D13").Value = '2.308.78'
$ws.Range("E13").Value = '  -0.31%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.15'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.77%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.801'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -6.38%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.25'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.65%  '

$ws.Range("D17").Value = '2.003.07'
$ws.Range("E17").Value = '  -1.45%  '

$ws.Range("D18").Value = '37.300.66'
$ws.Range("E18").Value = '  +0.54%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '70.17'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.64%  '

$ws.Range("D20").Value = '0.0₃0840'
$ws.Range("E20").Value = '  -3.16%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '234.05'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.42%  '

$ws.Range("E22").Value = '  -2.01%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.62'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.18%  '

$ws.Range("E24").Value = '  +0.01%  '

$ws.Range("E25").Value = '  -0.67%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.98'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.74%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.92'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.80%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.68'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.65%  '

$ws.Range("E29").Value = '  -5.31%  '

$ws.Range("E30").Value = '  -2.36%  '

$ws.Range("E31").Value = '  -1.41%  '

$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0646'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.51%  '

$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.61'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.63%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.54'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.34%  '

$ws.Range("E35").Value = '  -3.99%  '

$ws.Range("E36").Value = '  +0.46%  '

$ws.Range("E37").Value = '  -0.13%  '

$ws.Range("E38").Value = '  -3.65%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.37'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.37%  '

$ws.Range("E40").Value = '  +3.88%  '

$ws.Range("E41").Value = '  +0.54%  '

$ws.Range("E42").Value = '  -0.90%  '

$ws.Range("E43").Value = '  -5.67%  '

$ws.Range("D44").Value = '1.417.56'
$ws.Range("E44").Value = '  +1.75%  '

$ws.Range("E45").Value = '  -5.32%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '90.05'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.29%  '

$ws.Range("E47").Value = '  -2.99%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.07'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.18%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.93'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.56%  '

$ws.Range("D50").Value = '2.200.41'
$ws.Range("E50").Value = '  -0.35%  '

$ws.Range("E51").Value = '  -10.99%  '
